# Atualizando o arquivo XLSX
# Apply updated odds values for row 7 (match h0fBkCa6)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G7").Value = 1.5
$ws.Range("H7").Value = 3.6
$ws.Range("I7").Value = 7.5
$ws.Range("L7").Value = 8
$ws.Range("U7").Value = 1.57
$ws.Range("V7").Value = 2.25
$ws.Range("W7").Value = 2.75
$ws.Range("X7").Value = 1.4
$ws.Range("Y7").Value = 4.5
$ws.Range("AA7").Value = 10
$ws.Range("AE7").Value = 6.5
$ws.Range("AH7").Value = 151
$ws.Range("AK7").Value = 34
$ws.Range("AL7").Value = 23
$ws.Range("AN7").Value = 67
$ws.Range("AO7").Value = 81
$ws.Range("AR7").Value = 4
$ws.Range("AS7").Value = 1.23
